# Update header labels from "OutSample" to "TestSample"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "RMSE_TestSample"
$ws.Range("F1").Value = "R2_TestSample"
$ws.Range("G1").Value = "Adjusted_R2_TestSample"

# Update the E/F/G data values (rows 2-9) with the corrected
# test-sample metrics from the proper train-test-validate split.
$values = @{
    2 = @(0.1801641115658563, 0.9454294466443746, 0.9403766176299648)
    3 = @(0.1736332819363075, 0.9510423763266515, 0.945500003835329)
    4 = @(0.1632675558436219, 0.9556315991925847, 0.9496589298531249)
    5 = @(0.1542706850041704, 0.9625584272765664, 0.9566852393983807)
    6 = @(0.1469574528607001, 0.9674561806886989, 0.9615982932126647)
    7 = @(0.1391661125983226, 0.9712112835033809, 0.9653360352387648)
    8 = @(0.1318358622157212, 0.9743259144451667, 0.9684422698388507)
    9 = @(0.1266954201762954, 0.976879555937572, 0.9709764638365266)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Range("E$row").Value = $rowValues[0]
    $ws.Range("F$row").Value = $rowValues[1]
    $ws.Range("G$row").Value = $rowValues[2]
}
